$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The job postings beyond the first five (rows 7 through 16) are being
# removed, keeping only the header row and the first five job rows.
$ws.Range("A7:D16").EntireRow.Delete()

# Capture the hyperlink targets that must survive (rows 2-6 in column D)
# before clearing the worksheet's hyperlink collection, since deleting
# rows above does not automatically drop the now-orphaned hyperlink
# definitions that used to live on rows 7-16.
$keepTargets = @()
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $link = $ws.Hyperlinks.Item($i)
    if ($link.Range.Row -le 6) {
        $keepTargets += , @($link.Range.Row, $link.Range.Text)
    }
}

# Remove every hyperlink definition on the sheet (this also clears the
# stale ones that pointed at the now-deleted rows 7-16).
$ws.Hyperlinks.Delete()

# Re-create only the hyperlinks that belong to the surviving rows 2-6,
# then restore the built-in "Hyperlink" cell style so the cell formatting
# matches what it was originally.
foreach ($entry in $keepTargets) {
    $row = $entry[0]
    $target = $entry[1]
    $cell = $ws.Cells.Item($row, 4)
    $ws.Hyperlinks.Add($cell, $target)
    $cell.Style = "Hyperlink"
}
